$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 58
$ws.Range("B58").Value = 2632502
$ws.Range("F58").Value = 'GKS Belchatow'
$ws.Range("G58").Value = 'Stomil Olsztyn'
$ws.Range("H58").Value = 2
$ws.Range("I58").Value = 3
$ws.Range("J58").Value = 'A'
$ws.Range("K58").Value = 2.15
$ws.Range("L58").Value = 3.2
$ws.Range("M58").Value = 3.1
$ws.Range("N58").Value = 2.4
$ws.Range("O58").Value = 3.2
$ws.Range("P58").Value = 3
$ws.Range("Q58").Value = 0
$ws.Range("R58").Value = 1.725
$ws.Range("S58").Value = 2.075
$ws.Range("T58").Value = 1.75
$ws.Range("U58").Value = 1.775
$ws.Range("V58").Value = 2.025
$ws.Range("W58").Value = -1
$ws.Range("X58").Value = -1
$ws.Range("Y58").Value = 2
$ws.Range("Z58").Value = -1
$ws.Range("AA58").Value = 1.075
$ws.Range("AB58").Value = 0.7749999999999999
$ws.Range("AC58").Value = -1

# Row 60
$ws.Range("B60").Value = 2632496
$ws.Range("F60").Value = 'GKS Jastrzebie'
$ws.Range("G60").Value = 'Gornik Leczna'
$ws.Range("H60").Value = 1
$ws.Range("I60").Value = 1
$ws.Range("J60").Value = 'D'
$ws.Range("K60").Value = 3.5
$ws.Range("L60").Value = 3.25
$ws.Range("M60").Value = 1.95
$ws.Range("N60").Value = 3
$ws.Range("O60").Value = 2.9
$ws.Range("P60").Value = 2.3
$ws.Range("Q60").Value = 0.25
$ws.Range("R60").Value = 1.775
$ws.Range("S60").Value = 2.025
$ws.Range("T60").Value = 2
$ws.Range("U60").Value = 1.775
$ws.Range("V60").Value = 2.025
$ws.Range("W60").Value = -1
$ws.Range("X60").Value = 1.9
$ws.Range("Y60").Value = -1
$ws.Range("Z60").Value = 0.3875
$ws.Range("AA60").Value = -0.5
$ws.Range("AB60").Value = 0
$ws.Range("AC60").Value = -0

# Row 61
$ws.Range("B61").Value = 2632497
$ws.Range("F61").Value = 'Chrobry Glogow'
$ws.Range("G61").Value = 'Miedz Legnica'
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 'D'
$ws.Range("K61").Value = 2.875
$ws.Range("L61").Value = 3.1
$ws.Range("M61").Value = 2.3
$ws.Range("N61").Value = 2.9
$ws.Range("O61").Value = 3.1
$ws.Range("P61").Value = 2.25
$ws.Range("Q61").Value = 0.25
$ws.Range("R61").Value = 1.775
$ws.Range("S61").Value = 2.025
$ws.Range("T61").Value = 2.25
$ws.Range("U61").Value = 2
$ws.Range("V61").Value = 1.8
$ws.Range("W61").Value = -1
$ws.Range("X61").Value = 2.1
$ws.Range("Y61").Value = -1
$ws.Range("Z61").Value = 0.3875
$ws.Range("AA61").Value = -0.5
$ws.Range("AB61").Value = -1
$ws.Range("AC61").Value = 0.8

# Row 62
$ws.Range("B62").Value = 2632500
$ws.Range("F62").Value = 'Arka Gdynia'
$ws.Range("G62").Value = 'Odra Opole'
$ws.Range("H62").Value = 1
$ws.Range("I62").Value = 1
$ws.Range("J62").Value = 'D'
$ws.Range("K62").Value = 1.6
$ws.Range("L62").Value = 4
$ws.Range("M62").Value = 4.5
$ws.Range("N62").Value = 1.727
$ws.Range("O62").Value = 3.4
$ws.Range("P62").Value = 4.75
$ws.Range("Q62").Value = -0.5
$ws.Range("R62").Value = 1.8
$ws.Range("S62").Value = 2.05
$ws.Range("T62").Value = 2
$ws.Range("U62").Value = 1.8
$ws.Range("V62").Value = 2.05
$ws.Range("W62").Value = -1
$ws.Range("X62").Value = 2.4
$ws.Range("Y62").Value = -1
$ws.Range("Z62").Value = -1
$ws.Range("AA62").Value = 1.05
$ws.Range("AB62").Value = 0
$ws.Range("AC62").Value = -0

# Row 82
$ws.Range("B82").Value = 2632523
$ws.Range("F82").Value = 'Sandecja Nowy Sacz'
$ws.Range("G82").Value = 'Resovia Rzeszow'
$ws.Range("H82").Value = 1
$ws.Range("I82").Value = 1
$ws.Range("J82").Value = 'D'
$ws.Range("K82").Value = 1.833
$ws.Range("L82").Value = 3.5
$ws.Range("M82").Value = 3.8
$ws.Range("N82").Value = 1.7
$ws.Range("O82").Value = 3.5
$ws.Range("P82").Value = 4.5
$ws.Range("Q82").Value = -0.75
$ws.Range("R82").Value = 2.1
$ws.Range("S82").Value = 1.775
$ws.Range("T82").Value = 2.25
$ws.Range("U82").Value = 2.1
$ws.Range("V82").Value = 1.775
$ws.Range("W82").Value = -1
$ws.Range("X82").Value = 2.5
$ws.Range("Y82").Value = -1
$ws.Range("Z82").Value = -1
$ws.Range("AA82").Value = 0.7749999999999999
$ws.Range("AB82").Value = -0.5
$ws.Range("AC82").Value = 0.3875

# Row 83
$ws.Range("B83").Value = 2632525
$ws.Range("F83").Value = 'Stomil Olsztyn'
$ws.Range("G83").Value = 'Korona Kielce'
$ws.Range("H83").Value = 2
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 'H'
$ws.Range("K83").Value = 2.75
$ws.Range("L83").Value = 3.2
$ws.Range("M83").Value = 2.4
$ws.Range("N83").Value = 3.5
$ws.Range("O83").Value = 3.2
$ws.Range("P83").Value = 2.05
$ws.Range("Q83").Value = 0.25
$ws.Range("R83").Value = 2
$ws.Range("S83").Value = 1.85
$ws.Range("T83").Value = 2
$ws.Range("U83").Value = 1.975
$ws.Range("V83").Value = 1.875
$ws.Range("W83").Value = 2.5
$ws.Range("X83").Value = -1
$ws.Range("Y83").Value = -1
$ws.Range("Z83").Value = 1
$ws.Range("AA83").Value = -1
$ws.Range("AB83").Value = 0
$ws.Range("AC83").Value = -0

# Row 93
$ws.Range("B93").Value = 2632534
$ws.Range("F93").Value = 'GKS Jastrzebie'
$ws.Range("G93").Value = 'Odra Opole'
$ws.Range("H93").Value = 3
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 'H'
$ws.Range("K93").Value = 2.625
$ws.Range("L93").Value = 3.2
$ws.Range("M93").Value = 2.5
$ws.Range("N93").Value = 2.4
$ws.Range("O93").Value = 3.2
$ws.Range("P93").Value = 2.75
$ws.Range("Q93").Value = 0
$ws.Range("R93").Value = 1.775
$ws.Range("S93").Value = 2.025
$ws.Range("T93").Value = 2
$ws.Range("U93").Value = 1.975
$ws.Range("V93").Value = 1.825
$ws.Range("W93").Value = 1.4
$ws.Range("X93").Value = -1
$ws.Range("Y93").Value = -1
$ws.Range("Z93").Value = 0.7749999999999999
$ws.Range("AA93").Value = -1
$ws.Range("AB93").Value = 0.9750000000000001
$ws.Range("AC93").Value = -1

# Row 94
$ws.Range("B94").Value = 2632532
$ws.Range("F94").Value = 'Gornik Leczna'
$ws.Range("G94").Value = 'Termalica BB Nieciecza'
$ws.Range("H94").Value = 3
$ws.Range("I94").Value = 3
$ws.Range("J94").Value = 'D'
$ws.Range("K94").Value = 2.625
$ws.Range("L94").Value = 3.2
$ws.Range("M94").Value = 2.5
$ws.Range("N94").Value = 3.4
$ws.Range("O94").Value = 3.2
$ws.Range("P94").Value = 2.2
$ws.Range("Q94").Value = 0.25
$ws.Range("R94").Value = 1.825
$ws.Range("S94").Value = 1.975
$ws.Range("T94").Value = 2.25
$ws.Range("U94").Value = 1.975
$ws.Range("V94").Value = 1.825
$ws.Range("W94").Value = -1
$ws.Range("X94").Value = 2.2
$ws.Range("Y94").Value = -1
$ws.Range("Z94").Value = 0.4125
$ws.Range("AA94").Value = -0.5
$ws.Range("AB94").Value = 0.9750000000000001
$ws.Range("AC94").Value = -1

# Row 148
$ws.Range("B148").Value = 2755807
$ws.Range("F148").Value = 'Miedz Legnica'
$ws.Range("G148").Value = 'Zaglebie Sosnowiec'
$ws.Range("H148").Value = 1
$ws.Range("I148").Value = 1
$ws.Range("J148").Value = 'D'
$ws.Range("K148").Value = 1.95
$ws.Range("L148").Value = 3.6
$ws.Range("M148").Value = 3.2
$ws.Range("N148").Value = 1.727
$ws.Range("O148").Value = 4
$ws.Range("P148").Value = 4
$ws.Range("Q148").Value = -0.75
$ws.Range("R148").Value = 2.025
$ws.Range("S148").Value = 1.825
$ws.Range("T148").Value = 2.75
$ws.Range("U148").Value = 2
$ws.Range("V148").Value = 1.85
$ws.Range("W148").Value = -1
$ws.Range("X148").Value = 3
$ws.Range("Y148").Value = -1
$ws.Range("Z148").Value = -1
$ws.Range("AA148").Value = 0.825
$ws.Range("AB148").Value = -1
$ws.Range("AC148").Value = 0.8500000000000001

# Row 149
$ws.Range("B149").Value = 2755798
$ws.Range("F149").Value = 'Chrobry Glogow'
$ws.Range("G149").Value = 'Arka Gdynia'
$ws.Range("H149").Value = 0
$ws.Range("I149").Value = 3
$ws.Range("J149").Value = 'A'
$ws.Range("K149").Value = 3.75
$ws.Range("L149").Value = 3.3
$ws.Range("M149").Value = 2
$ws.Range("N149").Value = 3.7
$ws.Range("O149").Value = 3.3
$ws.Range("P149").Value = 2
$ws.Range("Q149").Value = 0.25
$ws.Range("R149").Value = 2
$ws.Range("S149").Value = 1.8
$ws.Range("T149").Value = 2.5
$ws.Range("U149").Value = 1.975
$ws.Range("V149").Value = 1.825
$ws.Range("W149").Value = -1
$ws.Range("X149").Value = -1
$ws.Range("Y149").Value = 1
$ws.Range("Z149").Value = -1
$ws.Range("AA149").Value = 0.8
$ws.Range("AB149").Value = 0.9750000000000001
$ws.Range("AC149").Value = -1

# Row 150
$ws.Range("B150").Value = 2759298
$ws.Range("F150").Value = 'Gornik Leczna'
$ws.Range("G150").Value = 'Sandecja Nowy Sacz'
$ws.Range("H150").Value = 3
$ws.Range("I150").Value = 0
$ws.Range("J150").Value = 'H'
$ws.Range("K150").Value = 2
$ws.Range("L150").Value = 3.6
$ws.Range("M150").Value = 3.2
$ws.Range("N150").Value = 2.1
$ws.Range("O150").Value = 3.6
$ws.Range("P150").Value = 3.2
$ws.Range("Q150").Value = -0.25
$ws.Range("R150").Value = 1.85
$ws.Range("S150").Value = 1.95
$ws.Range("T150").Value = 2.5
$ws.Range("U150").Value = 1.825
$ws.Range("V150").Value = 1.975
$ws.Range("W150").Value = 1.1
$ws.Range("X150").Value = -1
$ws.Range("Y150").Value = -1
$ws.Range("Z150").Value = 0.8500000000000001
$ws.Range("AA150").Value = -1
$ws.Range("AB150").Value = 0.825
$ws.Range("AC150").Value = -1

# Row 151
$ws.Range("B151").Value = 3611527
$ws.Range("F151").Value = 'GKS Jastrzebie'
$ws.Range("G151").Value = 'Resovia Rzeszow'
$ws.Range("H151").Value = 2
$ws.Range("I151").Value = 1
$ws.Range("J151").Value = 'H'
$ws.Range("K151").Value = 2.45
$ws.Range("L151").Value = 3.4
$ws.Range("M151").Value = 2.55
$ws.Range("N151").Value = 2.625
$ws.Range("O151").Value = 3.3
$ws.Range("P151").Value = 2.625
$ws.Range("Q151").Value = 0
$ws.Range("R151").Value = 1.925
$ws.Range("S151").Value = 1.875
$ws.Range("T151").Value = 2.25
$ws.Range("U151").Value = 1.775
$ws.Range("V151").Value = 2.025
$ws.Range("W151").Value = 1.625
$ws.Range("X151").Value = -1
$ws.Range("Y151").Value = -1
$ws.Range("Z151").Value = 0.925
$ws.Range("AA151").Value = -1
$ws.Range("AB151").Value = 0.7749999999999999
$ws.Range("AC151").Value = -1

# Row 152
$ws.Range("B152").Value = 2758744
$ws.Range("F152").Value = 'Widzew Lodz'
$ws.Range("G152").Value = 'Odra Opole'
$ws.Range("H152").Value = 2
$ws.Range("I152").Value = 1
$ws.Range("J152").Value = 'H'
$ws.Range("K152").Value = 1.8
$ws.Range("L152").Value = 3.6
$ws.Range("M152").Value = 4.2
$ws.Range("N152").Value = 1.8
$ws.Range("O152").Value = 3.6
$ws.Range("P152").Value = 4.2
$ws.Range("Q152").Value = -0.5
$ws.Range("R152").Value = 1.8
$ws.Range("S152").Value = 2
$ws.Range("T152").Value = 2.5
$ws.Range("U152").Value = 1.95
$ws.Range("V152").Value = 1.85
$ws.Range("W152").Value = 0.8
$ws.Range("X152").Value = -1
$ws.Range("Y152").Value = -1
$ws.Range("Z152").Value = 0.8
$ws.Range("AA152").Value = -1
$ws.Range("AB152").Value = 0.95
$ws.Range("AC152").Value = -1

# Row 153
$ws.Range("B153").Value = 2759740
$ws.Range("F153").Value = 'MKS Puszcza Niepolomice'
$ws.Range("G153").Value = 'GKS Belchatow'
$ws.Range("H153").Value = 2
$ws.Range("I153").Value = 0
$ws.Range("J153").Value = 'H'
$ws.Range("K153").Value = 2.3
$ws.Range("L153").Value = 3.4
$ws.Range("M153").Value = 3
$ws.Range("N153").Value = 2.3
$ws.Range("O153").Value = 3.4
$ws.Range("P153").Value = 3
$ws.Range("Q153").Value = -0.25
$ws.Range("R153").Value = 2.025
$ws.Range("S153").Value = 1.775
$ws.Range("T153").Value = 2.5
$ws.Range("U153").Value = 1.975
$ws.Range("V153").Value = 1.825
$ws.Range("W153").Value = 1.3
$ws.Range("X153").Value = -1
$ws.Range("Y153").Value = -1
$ws.Range("Z153").Value = 1.025
$ws.Range("AA153").Value = -1
$ws.Range("AB153").Value = -1
$ws.Range("AC153").Value = 0.825

# Row 154
$ws.Range("B154").Value = 2759543
$ws.Range("F154").Value = 'GKS Tychy 71'
$ws.Range("G154").Value = 'LKS Lodz'
$ws.Range("H154").Value = 1
$ws.Range("I154").Value = 1
$ws.Range("J154").Value = 'D'
$ws.Range("K154").Value = 2.375
$ws.Range("L154").Value = 3.3
$ws.Range("M154").Value = 2.7
$ws.Range("N154").Value = 1.8
$ws.Range("O154").Value = 4
$ws.Range("P154").Value = 3.8
$ws.Range("Q154").Value = -0.5
$ws.Range("R154").Value = 1.8
$ws.Range("S154").Value = 2
$ws.Range("T154").Value = 2.5
$ws.Range("U154").Value = 1.825
$ws.Range("V154").Value = 1.975
$ws.Range("W154").Value = -1
$ws.Range("X154").Value = 3
$ws.Range("Y154").Value = -1
$ws.Range("Z154").Value = -1
$ws.Range("AA154").Value = 1
$ws.Range("AB154").Value = -1
$ws.Range("AC154").Value = 0.9750000000000001

# Row 155
$ws.Range("B155").Value = 2758703
$ws.Range("F155").Value = 'Termalica BB Nieciecza'
$ws.Range("G155").Value = 'Stomil Olsztyn'
$ws.Range("H155").Value = 0
$ws.Range("I155").Value = 0
$ws.Range("J155").Value = 'D'
$ws.Range("K155").Value = 1.285
$ws.Range("L155").Value = 5
$ws.Range("M155").Value = 9
$ws.Range("N155").Value = 1.3
$ws.Range("O155").Value = 5
$ws.Range("P155").Value = 10
$ws.Range("Q155").Value = -1.5
$ws.Range("R155").Value = 1.925
$ws.Range("S155").Value = 1.925
$ws.Range("T155").Value = 2.75
$ws.Range("U155").Value = 1.775
$ws.Range("V155").Value = 2.025
$ws.Range("W155").Value = -1
$ws.Range("X155").Value = 4
$ws.Range("Y155").Value = -1
$ws.Range("Z155").Value = -1
$ws.Range("AA155").Value = 0.925
$ws.Range("AB155").Value = -1
$ws.Range("AC155").Value = 1.025

# Row 156
$ws.Range("B156").Value = 2761522
$ws.Range("F156").Value = 'Radomiak Radom'
$ws.Range("G156").Value = 'Korona Kielce'
$ws.Range("H156").Value = 2
$ws.Range("I156").Value = 0
$ws.Range("J156").Value = 'H'
$ws.Range("K156").Value = 1.45
$ws.Range("L156").Value = 4.2
$ws.Range("M156").Value = 7
$ws.Range("N156").Value = 1.45
$ws.Range("O156").Value = 4.2
$ws.Range("P156").Value = 7
$ws.Range("Q156").Value = -1
$ws.Range("R156").Value = 1.775
$ws.Range("S156").Value = 2.025
$ws.Range("T156").Value = 2.5
$ws.Range("U156").Value = 1.9
$ws.Range("V156").Value = 1.9
$ws.Range("W156").Value = 0.45
$ws.Range("X156").Value = -1
$ws.Range("Y156").Value = -1
$ws.Range("Z156").Value = 0.7749999999999999
$ws.Range("AA156").Value = -1
$ws.Range("AB156").Value = -1
$ws.Range("AC156").Value = 0.8999999999999999

# Row 187
$ws.Range("B187").Value = 3724141
$ws.Range("F187").Value = 'Miedz Legnica'
$ws.Range("G187").Value = 'Skra Czestochowa'
$ws.Range("H187").Value = 1
$ws.Range("I187").Value = 1
$ws.Range("J187").Value = 'D'
$ws.Range("K187").Value = 1.444
$ws.Range("L187").Value = 4.5
$ws.Range("M187").Value = 6.5
$ws.Range("N187").Value = 1.333
$ws.Range("O187").Value = 5.25
$ws.Range("P187").Value = 8
$ws.Range("Q187").Value = -1.25
$ws.Range("R187").Value = 1.825
$ws.Range("S187").Value = 1.975
$ws.Range("T187").Value = 2.75
$ws.Range("U187").Value = 1.925
$ws.Range("V187").Value = 1.875
$ws.Range("W187").Value = -1
$ws.Range("X187").Value = 4.25
$ws.Range("Y187").Value = -1
$ws.Range("Z187").Value = -1
$ws.Range("AA187").Value = 0.9750000000000001
$ws.Range("AB187").Value = -1
$ws.Range("AC187").Value = 0.875

# Row 188
$ws.Range("B188").Value = 3723191
$ws.Range("F188").Value = 'Resovia Rzeszow'
$ws.Range("G188").Value = 'MKS Puszcza Niepolomice'
$ws.Range("H188").Value = 1
$ws.Range("I188").Value = 1
$ws.Range("J188").Value = 'D'
$ws.Range("K188").Value = 1.75
$ws.Range("L188").Value = 3.6
$ws.Range("M188").Value = 4.5
$ws.Range("N188").Value = 1.909
$ws.Range("O188").Value = 3.6
$ws.Range("P188").Value = 4
$ws.Range("Q188").Value = -0.5
$ws.Range("R188").Value = 1.875
$ws.Range("S188").Value = 1.925
$ws.Range("T188").Value = 2.25
$ws.Range("U188").Value = 1.925
$ws.Range("V188").Value = 1.875
$ws.Range("W188").Value = -1
$ws.Range("X188").Value = 2.6
$ws.Range("Y188").Value = -1
$ws.Range("Z188").Value = -1
$ws.Range("AA188").Value = 0.925
$ws.Range("AB188").Value = -0.5
$ws.Range("AC188").Value = 0.4375

# Row 372
$ws.Range("B372").Value = 4782190
$ws.Range("F372").Value = 'Stomil Olsztyn'
$ws.Range("G372").Value = 'Gornik Polkowice'
$ws.Range("H372").Value = 0
$ws.Range("I372").Value = 4
$ws.Range("J372").Value = 'A'
$ws.Range("K372").Value = 2.25
$ws.Range("L372").Value = 3.2
$ws.Range("M372").Value = 3
$ws.Range("N372").Value = 2.2
$ws.Range("O372").Value = 3.1
$ws.Range("P372").Value = 3.2
$ws.Range("Q372").Value = -0.25
$ws.Range("R372").Value = 1.95
$ws.Range("S372").Value = 1.85
$ws.Range("T372").Value = 2.25
$ws.Range("U372").Value = 1.825
$ws.Range("V372").Value = 1.975
$ws.Range("W372").Value = -1
$ws.Range("X372").Value = -1
$ws.Range("Y372").Value = 2.2
$ws.Range("Z372").Value = -1
$ws.Range("AA372").Value = 0.8500000000000001
$ws.Range("AB372").Value = 0.825
$ws.Range("AC372").Value = -1

# Row 373
$ws.Range("B373").Value = 3724248
$ws.Range("F373").Value = 'Podbeskidzie Bielsko Biala'
$ws.Range("G373").Value = 'GKS Jastrzebie'
$ws.Range("H373").Value = 0
$ws.Range("I373").Value = 0
$ws.Range("J373").Value = 'D'
$ws.Range("K373").Value = 1.727
$ws.Range("L373").Value = 3.4
$ws.Range("M373").Value = 4.333
$ws.Range("N373").Value = 1.45
$ws.Range("O373").Value = 3.8
$ws.Range("P373").Value = 6
$ws.Range("Q373").Value = -1
$ws.Range("R373").Value = 1.8
$ws.Range("S373").Value = 2
$ws.Range("T373").Value = 2.5
$ws.Range("U373").Value = 1.95
$ws.Range("V373").Value = 1.85
$ws.Range("W373").Value = -1
$ws.Range("X373").Value = 2.8
$ws.Range("Y373").Value = -1
$ws.Range("Z373").Value = -1
$ws.Range("AA373").Value = 1
$ws.Range("AB373").Value = -1
$ws.Range("AC373").Value = 0.8500000000000001

# Row 456
$ws.Range("B456").Value = 3951941
$ws.Range("F456").Value = 'GKS Katowice'
$ws.Range("G456").Value = 'LKS Lodz'
$ws.Range("H456").Value = 2
$ws.Range("I456").Value = 0
$ws.Range("J456").Value = 'H'
$ws.Range("K456").Value = 2.25
$ws.Range("L456").Value = 3.4
$ws.Range("M456").Value = 2.9
$ws.Range("N456").Value = 2.2
$ws.Range("O456").Value = 3.4
$ws.Range("P456").Value = 3.25
$ws.Range("Q456").Value = -0.25
$ws.Range("R456").Value = 1.9
$ws.Range("S456").Value = 1.9
$ws.Range("T456").Value = 2.5
$ws.Range("U456").Value = 1.825
$ws.Range("V456").Value = 1.975
$ws.Range("W456").Value = 1.2
$ws.Range("X456").Value = -1
$ws.Range("Y456").Value = -1
$ws.Range("Z456").Value = 0.8999999999999999
$ws.Range("AA456").Value = -1
$ws.Range("AB456").Value = -1
$ws.Range("AC456").Value = 0.9750000000000001

# Row 457
$ws.Range("B457").Value = 3945948
$ws.Range("F457").Value = 'Stomil Olsztyn'
$ws.Range("G457").Value = 'MKS Puszcza Niepolomice'
$ws.Range("H457").Value = 1
$ws.Range("I457").Value = 0
$ws.Range("J457").Value = 'H'
$ws.Range("K457").Value = 2.875
$ws.Range("L457").Value = 3.4
$ws.Range("M457").Value = 2.25
$ws.Range("N457").Value = 2.75
$ws.Range("O457").Value = 3.5
$ws.Range("P457").Value = 2.3
$ws.Range("Q457").Value = 0.25
$ws.Range("R457").Value = 1.8
$ws.Range("S457").Value = 2.05
$ws.Range("T457").Value = 2.75
$ws.Range("U457").Value = 1.875
$ws.Range("V457").Value = 1.975
$ws.Range("W457").Value = 1.75
$ws.Range("X457").Value = -1
$ws.Range("Y457").Value = -1
$ws.Range("Z457").Value = 0.8
$ws.Range("AA457").Value = -1
$ws.Range("AB457").Value = -1
$ws.Range("AC457").Value = 0.9750000000000001

# Row 458
$ws.Range("B458").Value = 3956530
$ws.Range("F458").Value = 'Chrobry Glogow'
$ws.Range("G458").Value = 'Zaglebie Sosnowiec'
$ws.Range("H458").Value = 4
$ws.Range("I458").Value = 0
$ws.Range("J458").Value = 'H'
$ws.Range("K458").Value = 1.666
$ws.Range("L458").Value = 3.6
$ws.Range("M458").Value = 4.75
$ws.Range("N458").Value = 2.05
$ws.Range("O458").Value = 3.5
$ws.Range("P458").Value = 3.25
$ws.Range("Q458").Value = -0.5
$ws.Range("R458").Value = 2.025
$ws.Range("S458").Value = 1.775
$ws.Range("T458").Value = 2.75
$ws.Range("U458").Value = 1.9
$ws.Range("V458").Value = 1.9
$ws.Range("W458").Value = 1.05
$ws.Range("X458").Value = -1
$ws.Range("Y458").Value = -1
$ws.Range("Z458").Value = 1.025
$ws.Range("AA458").Value = -1
$ws.Range("AB458").Value = 0.8999999999999999
$ws.Range("AC458").Value = -1

# Row 459
$ws.Range("B459").Value = 3953056
$ws.Range("F459").Value = 'Widzew Lodz'
$ws.Range("G459").Value = 'Podbeskidzie Bielsko Biala'
$ws.Range("H459").Value = 2
$ws.Range("I459").Value = 1
$ws.Range("J459").Value = 'H'
$ws.Range("K459").Value = 2.05
$ws.Range("L459").Value = 3.4
$ws.Range("M459").Value = 3.3
$ws.Range("N459").Value = 2.05
$ws.Range("O459").Value = 3.6
$ws.Range("P459").Value = 3.1
$ws.Range("Q459").Value = -0.25
$ws.Range("R459").Value = 1.825
$ws.Range("S459").Value = 1.975
$ws.Range("T459").Value = 2.75
$ws.Range("U459").Value = 1.85
$ws.Range("V459").Value = 1.95
$ws.Range("W459").Value = 1.05
$ws.Range("X459").Value = -1
$ws.Range("Y459").Value = -1
$ws.Range("Z459").Value = 0.825
$ws.Range("AA459").Value = -1
$ws.Range("AB459").Value = 0.425
$ws.Range("AC459").Value = -0.5

# Row 462
$ws.Range("B462").Value = 3945949
$ws.Range("F462").Value = 'Gornik Polkowice'
$ws.Range("G462").Value = 'Miedz Legnica'
$ws.Range("H462").Value = 0
$ws.Range("I462").Value = 1
$ws.Range("J462").Value = 'A'
$ws.Range("K462").Value = 3.1
$ws.Range("L462").Value = 3.3
$ws.Range("M462").Value = 2.2
$ws.Range("N462").Value = 1.6
$ws.Range("O462").Value = 4
$ws.Range("P462").Value = 5.25
$ws.Range("Q462").Value = -0.75
$ws.Range("R462").Value = 1.775
$ws.Range("S462").Value = 2.025
$ws.Range("T462").Value = 2.75
$ws.Range("U462").Value = 1.975
$ws.Range("V462").Value = 1.825
$ws.Range("W462").Value = -1
$ws.Range("X462").Value = -1
$ws.Range("Y462").Value = 4.25
$ws.Range("Z462").Value = -1
$ws.Range("AA462").Value = 1.025
$ws.Range("AB462").Value = -1
$ws.Range("AC462").Value = 0.825

# Row 463
$ws.Range("B463").Value = 3956529
$ws.Range("F463").Value = 'Arka Gdynia'
$ws.Range("G463").Value = 'Sandecja Nowy Sacz'
$ws.Range("H463").Value = 2
$ws.Range("I463").Value = 1
$ws.Range("J463").Value = 'H'
$ws.Range("K463").Value = 1.666
$ws.Range("L463").Value = 3.75
$ws.Range("M463").Value = 4.5
$ws.Range("N463").Value = 1.727
$ws.Range("O463").Value = 3.8
$ws.Range("P463").Value = 4
$ws.Range("Q463").Value = -0.75
$ws.Range("R463").Value = 2
$ws.Range("S463").Value = 1.8
$ws.Range("T463").Value = 2.75
$ws.Range("U463").Value = 1.85
$ws.Range("V463").Value = 1.95
$ws.Range("W463").Value = 0.7270000000000001
$ws.Range("X463").Value = -1
$ws.Range("Y463").Value = -1
$ws.Range("Z463").Value = 0.5
$ws.Range("AA463").Value = -0.5
$ws.Range("AB463").Value = 0.425
$ws.Range("AC463").Value = -0.5

# Row 464
$ws.Range("B464").Value = 3951940
$ws.Range("F464").Value = 'GKS Tychy 71'
$ws.Range("G464").Value = 'Korona Kielce'
$ws.Range("H464").Value = 1
$ws.Range("I464").Value = 1
$ws.Range("J464").Value = 'D'
$ws.Range("K464").Value = 2.1
$ws.Range("L464").Value = 3.4
$ws.Range("M464").Value = 3.2
$ws.Range("N464").Value = 2.1
$ws.Range("O464").Value = 3.4
$ws.Range("P464").Value = 3.2
$ws.Range("Q464").Value = -0.25
$ws.Range("R464").Value = 1.8
$ws.Range("S464").Value = 2
$ws.Range("T464").Value = 2.75
$ws.Range("U464").Value = 2
$ws.Range("V464").Value = 1.8
$ws.Range("W464").Value = -1
$ws.Range("X464").Value = 2.4
$ws.Range("Y464").Value = -1
$ws.Range("Z464").Value = -0.5
$ws.Range("AA464").Value = 0.5
$ws.Range("AB464").Value = -1
$ws.Range("AC464").Value = 0.8

# Row 510
$ws.Range("B510").Value = 5138951
$ws.Range("F510").Value = 'Termalica BB Nieciecza'
$ws.Range("G510").Value = 'Chrobry Glogow'
$ws.Range("H510").Value = 3
$ws.Range("I510").Value = 1
$ws.Range("J510").Value = 'H'
$ws.Range("K510").Value = 1.666
$ws.Range("L510").Value = 3.6
$ws.Range("M510").Value = 4.333
$ws.Range("N510").Value = 1.666
$ws.Range("O510").Value = 3.5
$ws.Range("P510").Value = 4.5
$ws.Range("Q510").Value = -0.75
$ws.Range("R510").Value = 1.875
$ws.Range("S510").Value = 1.925
$ws.Range("T510").Value = 2.5
$ws.Range("U510").Value = 1.875
$ws.Range("V510").Value = 1.925
$ws.Range("W510").Value = 0.6659999999999999
$ws.Range("X510").Value = -1
$ws.Range("Y510").Value = -1
$ws.Range("Z510").Value = 0.875
$ws.Range("AA510").Value = -1
$ws.Range("AB510").Value = 0.875
$ws.Range("AC510").Value = -1

# Row 511
$ws.Range("B511").Value = 5140780
$ws.Range("F511").Value = 'Chojniczanka Chojnice'
$ws.Range("G511").Value = 'Stal Rzeszow'
$ws.Range("H511").Value = 1
$ws.Range("I511").Value = 1
$ws.Range("J511").Value = 'D'
$ws.Range("K511").Value = 2.55
$ws.Range("L511").Value = 3.2
$ws.Range("M511").Value = 2.55
$ws.Range("N511").Value = 2.625
$ws.Range("O511").Value = 3.4
$ws.Range("P511").Value = 2.375
$ws.Range("Q511").Value = 0
$ws.Range("R511").Value = 1.975
$ws.Range("S511").Value = 1.825
$ws.Range("T511").Value = 2.75
$ws.Range("U511").Value = 1.95
$ws.Range("V511").Value = 1.9
$ws.Range("W511").Value = -1
$ws.Range("X511").Value = 2.4
$ws.Range("Y511").Value = -1
$ws.Range("Z511").Value = 0
$ws.Range("AA511").Value = -0
$ws.Range("AB511").Value = -1
$ws.Range("AC511").Value = 0.8999999999999999

# Row 681
$ws.Range("B681").Value = 5139054
$ws.Range("F681").Value = 'GKS Tychy 71'
$ws.Range("G681").Value = 'Sandecja Nowy Sacz'
$ws.Range("H681").Value = 2
$ws.Range("I681").Value = 3
$ws.Range("J681").Value = 'A'
$ws.Range("K681").Value = 2.15
$ws.Range("L681").Value = 3.2
$ws.Range("M681").Value = 3.1
$ws.Range("N681").Value = 2.375
$ws.Range("O681").Value = 3
$ws.Range("P681").Value = 3
$ws.Range("Q681").Value = -0.25
$ws.Range("R681").Value = 2.025
$ws.Range("S681").Value = 1.775
$ws.Range("T681").Value = 2.25
$ws.Range("U681").Value = 1.975
$ws.Range("V681").Value = 1.825
$ws.Range("W681").Value = -1
$ws.Range("X681").Value = -1
$ws.Range("Y681").Value = 2
$ws.Range("Z681").Value = -1
$ws.Range("AA681").Value = 0.7749999999999999
$ws.Range("AB681").Value = 0.9750000000000001
$ws.Range("AC681").Value = -1

# Row 682
$ws.Range("B682").Value = 5139056
$ws.Range("F682").Value = 'Odra Opole'
$ws.Range("G682").Value = 'Arka Gdynia'
$ws.Range("H682").Value = 1
$ws.Range("I682").Value = 1
$ws.Range("J682").Value = 'D'
$ws.Range("K682").Value = 3.75
$ws.Range("L682").Value = 3.5
$ws.Range("M682").Value = 1.85
$ws.Range("N682").Value = 3.4
$ws.Range("O682").Value = 3.5
$ws.Range("P682").Value = 1.909
$ws.Range("Q682").Value = 0.5
$ws.Range("R682").Value = 1.85
$ws.Range("S682").Value = 2
$ws.Range("T682").Value = 2.75
$ws.Range("U682").Value = 2
$ws.Range("V682").Value = 1.85
$ws.Range("W682").Value = -1
$ws.Range("X682").Value = 2.5
$ws.Range("Y682").Value = -1
$ws.Range("Z682").Value = 0.8500000000000001
$ws.Range("AA682").Value = -1
$ws.Range("AB682").Value = -1
$ws.Range("AC682").Value = 0.8500000000000001

# Row 848
$ws.Range("B848").Value = 6803738
$ws.Range("F848").Value = 'Podbeskidzie Bielsko Biala'
$ws.Range("G848").Value = 'Gornik Leczna'
$ws.Range("H848").Value = 1
$ws.Range("I848").Value = 1
$ws.Range("J848").Value = 'D'
$ws.Range("K848").Value = 1.85
$ws.Range("L848").Value = 3.5
$ws.Range("M848").Value = 3.8
$ws.Range("N848").Value = 1.666
$ws.Range("O848").Value = 3.8
$ws.Range("P848").Value = 4.75
$ws.Range("Q848").Value = -0.75
$ws.Range("R848").Value = 1.825
$ws.Range("S848").Value = 1.975
$ws.Range("T848").Value = 2.5
$ws.Range("U848").Value = 1.825
$ws.Range("V848").Value = 1.975
$ws.Range("W848").Value = -1
$ws.Range("X848").Value = 2.8
$ws.Range("Y848").Value = -1
$ws.Range("Z848").Value = -1
$ws.Range("AA848").Value = 0.9750000000000001
$ws.Range("AB848").Value = -1
$ws.Range("AC848").Value = 0.9750000000000001

# Row 849
$ws.Range("B849").Value = 6803740
$ws.Range("F849").Value = 'Miedz Legnica'
$ws.Range("G849").Value = 'Odra Opole'
$ws.Range("H849").Value = 1
$ws.Range("I849").Value = 2
$ws.Range("J849").Value = 'A'
$ws.Range("K849").Value = 1.85
$ws.Range("L849").Value = 3.5
$ws.Range("M849").Value = 3.75
$ws.Range("N849").Value = 1.909
$ws.Range("O849").Value = 3.5
$ws.Range("P849").Value = 3.5
$ws.Range("Q849").Value = -0.5
$ws.Range("R849").Value = 1.975
$ws.Range("S849").Value = 1.825
$ws.Range("T849").Value = 2.25
$ws.Range("U849").Value = 1.9
$ws.Range("V849").Value = 1.9
$ws.Range("W849").Value = -1
$ws.Range("X849").Value = -1
$ws.Range("Y849").Value = 2.5
$ws.Range("Z849").Value = -1
$ws.Range("AA849").Value = 0.825
$ws.Range("AB849").Value = 0.8999999999999999
$ws.Range("AC849").Value = -1

# Row 930
$ws.Range("B930").Value = 6803793
$ws.Range("F930").Value = 'Odra Opole'
$ws.Range("G930").Value = 'Stal Rzeszow'
$ws.Range("H930").Value = 1
$ws.Range("I930").Value = 1
$ws.Range("J930").Value = 'D'
$ws.Range("K930").Value = 2.05
$ws.Range("L930").Value = 3.4
$ws.Range("M930").Value = 3.5
$ws.Range("N930").Value = 2.1
$ws.Range("O930").Value = 3.4
$ws.Range("P930").Value = 3.4
$ws.Range("Q930").Value = -0.25
$ws.Range("R930").Value = 1.825
$ws.Range("S930").Value = 1.975
$ws.Range("T930").Value = 2.5
$ws.Range("U930").Value = 1.95
$ws.Range("V930").Value = 1.85
$ws.Range("W930").Value = -1
$ws.Range("X930").Value = 2.4
$ws.Range("Y930").Value = -1
$ws.Range("Z930").Value = -0.5
$ws.Range("AA930").Value = 0.4875
$ws.Range("AB930").Value = -1
$ws.Range("AC930").Value = 0.8500000000000001

# Row 931
$ws.Range("B931").Value = 6803794
$ws.Range("F931").Value = 'Wisla Krakow'
$ws.Range("G931").Value = 'Gornik Leczna'
$ws.Range("H931").Value = 4
$ws.Range("I931").Value = 0
$ws.Range("J931").Value = 'H'
$ws.Range("K931").Value = 1.4
$ws.Range("L931").Value = 4.75
$ws.Range("M931").Value = 7
$ws.Range("N931").Value = 1.363
$ws.Range("O931").Value = 4.75
$ws.Range("P931").Value = 7.5
$ws.Range("Q931").Value = -1.25
$ws.Range("R931").Value = 1.8
$ws.Range("S931").Value = 2
$ws.Range("T931").Value = 2.75
$ws.Range("U931").Value = 1.775
$ws.Range("V931").Value = 2.025
$ws.Range("W931").Value = 0.363
$ws.Range("X931").Value = -1
$ws.Range("Y931").Value = -1
$ws.Range("Z931").Value = 0.8
$ws.Range("AA931").Value = -1
$ws.Range("AB931").Value = 0.7749999999999999
$ws.Range("AC931").Value = -1
